$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report data was regenerated: a couple of items were removed, a couple
# of new items were added, a few quantities were corrected, and the whole
# table ended up re-sorted/re-grouped by category. Net effect: the sheet now
# holds 22 data rows (rows 2-23) instead of 24 (rows 2-25).
# ---------------------------------------------------------------------------

# Drop the two rows that are no longer part of the table so the sheet shrinks
# from A1:E25 down to A1:E23.
$ws.Rows("24:25").Delete()

# Rewrite every data row (2-23) with the up to date report contents.
$ws.Cells.Item(2,1).Value = 26
$ws.Cells.Item(2,2).Value = "LIMPEZA"
$ws.Cells.Item(2,3).Value = "Sabaoliquido"
$ws.Cells.Item(2,4).Value = 600
$ws.Cells.Item(2,5).Value = "31/12/2025"

$ws.Cells.Item(3,1).Value = 26
$ws.Cells.Item(3,2).Value = "LIMPEZA"
$ws.Cells.Item(3,3).Value = "Sabaoembarra"
$ws.Cells.Item(3,4).Value = 2
$ws.Cells.Item(3,5).Value = ""

$ws.Cells.Item(4,1).Value = 26
$ws.Cells.Item(4,2).Value = "LIMPEZA"
$ws.Cells.Item(4,3).Value = "Borrachabranca"
$ws.Cells.Item(4,4).Value = 120
$ws.Cells.Item(4,5).Value = ""

$ws.Cells.Item(5,1).Value = 27
$ws.Cells.Item(5,2).Value = "FERRAMENTAS"
$ws.Cells.Item(5,3).Value = "Chavedefenda"
$ws.Cells.Item(5,4).Value = 1
$ws.Cells.Item(5,5).Value = ""

$ws.Cells.Item(6,1).Value = 20
$ws.Cells.Item(6,2).Value = "PAPELARIA"
$ws.Cells.Item(6,3).Value = "Cadernoespiral"
$ws.Cells.Item(6,4).Value = 50
$ws.Cells.Item(6,5).Value = ""

$ws.Cells.Item(7,1).Value = 20
$ws.Cells.Item(7,2).Value = "PAPELARIA"
$ws.Cells.Item(7,3).Value = "Borrachabranca"
$ws.Cells.Item(7,4).Value = 20
$ws.Cells.Item(7,5).Value = ""

$ws.Cells.Item(8,1).Value = 26
$ws.Cells.Item(8,2).Value = "LIMPEZA"
$ws.Cells.Item(8,3).Value = "Papela1Sulfite"
$ws.Cells.Item(8,4).Value = 20
$ws.Cells.Item(8,5).Value = ""

$ws.Cells.Item(9,1).Value = 26
$ws.Cells.Item(9,2).Value = "LIMPEZA"
$ws.Cells.Item(9,3).Value = "Detergente"
$ws.Cells.Item(9,4).Value = 1
$ws.Cells.Item(9,5).Value = ""

$ws.Cells.Item(10,1).Value = 26
$ws.Cells.Item(10,2).Value = "LIMPEZA"
$ws.Cells.Item(10,3).Value = "Cadernoespiral"
$ws.Cells.Item(10,4).Value = 299
$ws.Cells.Item(10,5).Value = ""

$ws.Cells.Item(11,1).Value = 29
$ws.Cells.Item(11,2).Value = "DIVERSOS"
$ws.Cells.Item(11,3).Value = "Garrafadeagua"
$ws.Cells.Item(11,4).Value = 20
# "11/06/2025" parses as a valid US-style date (Nov 6 2025), so Excel would
# otherwise silently convert it to a date serial. Force it to stay plain text
# the same way it is stored in the report, then strip the leftover number
# format so no extra cell style sticks around.
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "11/06/2025"
$ws.Cells.Item(11,5).ClearFormats()

$ws.Cells.Item(12,1).Value = 20
$ws.Cells.Item(12,2).Value = "PAPELARIA"
$ws.Cells.Item(12,3).Value = "Canetaesferograficaazul"
$ws.Cells.Item(12,4).Value = 300
$ws.Cells.Item(12,5).Value = "31/12/2025"

$ws.Cells.Item(13,1).Value = 24
$ws.Cells.Item(13,2).Value = "INFORMATICA"
$ws.Cells.Item(13,3).Value = "Mouseusb"
$ws.Cells.Item(13,4).Value = 1
$ws.Cells.Item(13,5).Value = ""

$ws.Cells.Item(14,1).Value = 28
$ws.Cells.Item(14,2).Value = "MATERIALESCOLAR"
$ws.Cells.Item(14,3).Value = "Cadernoespiral"
$ws.Cells.Item(14,4).Value = 100
$ws.Cells.Item(14,5).Value = ""

$ws.Cells.Item(15,1).Value = 24
$ws.Cells.Item(15,2).Value = "INFORMATICA"
$ws.Cells.Item(15,3).Value = "Cabovga"
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(15,5).Value = ""

$ws.Cells.Item(16,1).Value = 20
$ws.Cells.Item(16,2).Value = "PAPELARIA"
$ws.Cells.Item(16,3).Value = "Papela4Sulfite"
$ws.Cells.Item(16,4).Value = 15
$ws.Cells.Item(16,5).Value = ""

$ws.Cells.Item(17,1).Value = 29
$ws.Cells.Item(17,2).Value = "DIVERSOS"
$ws.Cells.Item(17,3).Value = "Acucar"
$ws.Cells.Item(17,4).Value = 50
$ws.Cells.Item(17,5).Value = "13/06/2025"

$ws.Cells.Item(18,1).Value = 20
$ws.Cells.Item(18,2).Value = "PAPELARIA"
$ws.Cells.Item(18,3).Value = "Canetapreta"
$ws.Cells.Item(18,4).Value = 4
# Same situation as row 11: "10/12/2025" parses as Oct 12 2025, so pin it as
# plain text and then clear the temporary format.
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "10/12/2025"
$ws.Cells.Item(18,5).ClearFormats()

$ws.Cells.Item(19,1).Value = 20
$ws.Cells.Item(19,2).Value = "PAPELARIA"
$ws.Cells.Item(19,3).Value = "Canetaazul"
$ws.Cells.Item(19,4).Value = 4
$ws.Cells.Item(19,5).Value = ""

$ws.Cells.Item(20,1).Value = 26
$ws.Cells.Item(20,2).Value = "LIMPEZA"
$ws.Cells.Item(20,3).Value = "Sabaoempo"
$ws.Cells.Item(20,4).Value = 2
$ws.Cells.Item(20,5).Value = ""

$ws.Cells.Item(21,1).Value = 24
$ws.Cells.Item(21,2).Value = "INFORMATICA"
$ws.Cells.Item(21,3).Value = "Tecladousb"
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = ""

$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "PAPELARIA"
$ws.Cells.Item(22,3).Value = "Canetavermelha"
$ws.Cells.Item(22,4).Value = 100
$ws.Cells.Item(22,5).Value = ""

$ws.Cells.Item(23,1).Value = 20
$ws.Cells.Item(23,2).Value = "PAPELARIA"
$ws.Cells.Item(23,3).Value = "Lapis"
$ws.Cells.Item(23,4).Value = 4
$ws.Cells.Item(23,5).Value = ""

# Narrow column B (Nome_Categoria) from 21 to 17 characters wide. ColumnWidth
# gets 5/6 of a character added internally when saved, so back that out here
# to land exactly on a stored width of 17.
$ws.Columns(2).ColumnWidth = 17 - (5/6)
